$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F,H:V data between the paired rows (58/59, 64/65, 124/125) ---
# (A-E — Indice/pais/torneio/temporada/data_partida — stay as-is for each row)

# Row 58
$ws.Cells.Item(58,6).Value = 'Las Palmas'
$ws.Cells.Item(58,7).Value = 1
$ws.Cells.Item(58,8).Value = 'Granada CF'
$ws.Cells.Item(58,9).Value = 0
$ws.Cells.Item(58,10).Value = 2.14
$ws.Cells.Item(58,11).Value = '11/09/2023 13:19'
$ws.Cells.Item(58,12).Value = 1.95
$ws.Cells.Item(58,13).Value = '24/09/2023 18:19'
$ws.Cells.Item(58,14).Value = 3.21
$ws.Cells.Item(58,15).Value = '11/09/2023 13:19'
$ws.Cells.Item(58,16).Value = 3.73
$ws.Cells.Item(58,17).Value = '24/09/2023 18:27'
$ws.Cells.Item(58,18).Value = 3.61
$ws.Cells.Item(58,19).Value = '11/09/2023 13:19'
$ws.Cells.Item(58,20).Value = 4.1
$ws.Cells.Item(58,21).Value = '24/09/2023 18:27'
$ws.Cells.Item(58,22).Value = 'https://www.betexplorer.com/football/spain/laliga/las-palmas-granada-cf/tWsBDE3N/'

# Row 59
$ws.Cells.Item(59,6).Value = 'Betis'
$ws.Cells.Item(59,7).Value = 1
$ws.Cells.Item(59,8).Value = 'Cadiz CF'
$ws.Cells.Item(59,9).Value = 1
$ws.Cells.Item(59,10).Value = 1.87
$ws.Cells.Item(59,11).Value = '05/09/2023 12:02'
$ws.Cells.Item(59,12).Value = 1.81
$ws.Cells.Item(59,13).Value = '24/09/2023 18:26'
$ws.Cells.Item(59,14).Value = 3.56
$ws.Cells.Item(59,15).Value = '05/09/2023 12:02'
$ws.Cells.Item(59,16).Value = 3.66
$ws.Cells.Item(59,17).Value = '24/09/2023 18:26'
$ws.Cells.Item(59,18).Value = 4.54
$ws.Cells.Item(59,19).Value = '05/09/2023 12:02'
$ws.Cells.Item(59,20).Value = 5.07
$ws.Cells.Item(59,21).Value = '24/09/2023 18:26'
$ws.Cells.Item(59,22).Value = 'https://www.betexplorer.com/football/spain/laliga/betis-cadiz/IicoJIZo/'

# Row 64
$ws.Cells.Item(64,6).Value = 'Ath Bilbao'
$ws.Cells.Item(64,7).Value = 2
$ws.Cells.Item(64,8).Value = 'Getafe'
$ws.Cells.Item(64,9).Value = 2
$ws.Cells.Item(64,10).Value = 1.71
$ws.Cells.Item(64,11).Value = '17/09/2023 09:02'
$ws.Cells.Item(64,12).Value = 1.53
$ws.Cells.Item(64,13).Value = '27/09/2023 18:31'
$ws.Cells.Item(64,14).Value = 3.42
$ws.Cells.Item(64,15).Value = '17/09/2023 09:02'
$ws.Cells.Item(64,16).Value = 4.06
$ws.Cells.Item(64,17).Value = '27/09/2023 18:49'
$ws.Cells.Item(64,18).Value = 5.44
$ws.Cells.Item(64,19).Value = '17/09/2023 09:02'
$ws.Cells.Item(64,20).Value = 7.73
$ws.Cells.Item(64,21).Value = '27/09/2023 18:49'
$ws.Cells.Item(64,22).Value = 'https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/'

# Row 65
$ws.Cells.Item(65,6).Value = 'Real Madrid'
$ws.Cells.Item(65,7).Value = 2
$ws.Cells.Item(65,8).Value = 'Las Palmas'
$ws.Cells.Item(65,9).Value = 0
$ws.Cells.Item(65,10).Value = 1.2
$ws.Cells.Item(65,11).Value = '23/09/2023 09:28'
$ws.Cells.Item(65,12).Value = 1.18
$ws.Cells.Item(65,13).Value = '27/09/2023 18:29'
$ws.Cells.Item(65,14).Value = 6.76
$ws.Cells.Item(65,15).Value = '23/09/2023 09:28'
$ws.Cells.Item(65,16).Value = 8
$ws.Cells.Item(65,17).Value = '27/09/2023 18:29'
$ws.Cells.Item(65,18).Value = 11.3
$ws.Cells.Item(65,19).Value = '23/09/2023 09:28'
$ws.Cells.Item(65,20).Value = 16.5
$ws.Cells.Item(65,21).Value = '27/09/2023 18:29'
$ws.Cells.Item(65,22).Value = 'https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/'

# Row 124
$ws.Cells.Item(124,6).Value = 'Granada CF'
$ws.Cells.Item(124,7).Value = 1
$ws.Cells.Item(124,8).Value = 'Getafe'
$ws.Cells.Item(124,9).Value = 1
$ws.Cells.Item(124,10).Value = 2.52
$ws.Cells.Item(124,11).Value = '29/10/2023 11:02'
$ws.Cells.Item(124,12).Value = 2.49
$ws.Cells.Item(124,13).Value = '11/11/2023 18:27'
$ws.Cells.Item(124,14).Value = 3.01
$ws.Cells.Item(124,15).Value = '29/10/2023 11:02'
$ws.Cells.Item(124,16).Value = 3.19
$ws.Cells.Item(124,17).Value = '11/11/2023 18:23'
$ws.Cells.Item(124,18).Value = 3.22
$ws.Cells.Item(124,19).Value = '29/10/2023 11:02'
$ws.Cells.Item(124,20).Value = 3.19
$ws.Cells.Item(124,21).Value = '11/11/2023 18:27'
$ws.Cells.Item(124,22).Value = 'https://www.betexplorer.com/football/spain/laliga/granada-cf-getafe/OSl3Qfr5/'

# Row 125
$ws.Cells.Item(125,6).Value = 'Osasuna'
$ws.Cells.Item(125,7).Value = 1
$ws.Cells.Item(125,8).Value = 'Las Palmas'
$ws.Cells.Item(125,9).Value = 1
$ws.Cells.Item(125,10).Value = 1.79
$ws.Cells.Item(125,11).Value = '29/10/2023 11:02'
$ws.Cells.Item(125,12).Value = 1.97
$ws.Cells.Item(125,13).Value = '11/11/2023 18:12'
$ws.Cells.Item(125,14).Value = 3.56
$ws.Cells.Item(125,15).Value = '29/10/2023 11:02'
$ws.Cells.Item(125,16).Value = 3.36
$ws.Cells.Item(125,17).Value = '11/11/2023 18:28'
$ws.Cells.Item(125,18).Value = 4.92
$ws.Cells.Item(125,19).Value = '29/10/2023 11:02'
$ws.Cells.Item(125,20).Value = 4.47
$ws.Cells.Item(125,21).Value = '11/11/2023 18:28'
$ws.Cells.Item(125,22).Value = 'https://www.betexplorer.com/football/spain/laliga/osasuna-las-palmas/UPiBOYCH/'

# --- Append new rows 130-137 ---

# Row 130
$ws.Range("A129:V129").Copy()
$ws.Range("A130:V130").PasteSpecial(-4122)
$ws.Cells.Item(130,1).Value = 129
$ws.Cells.Item(130,2).Value = 'spain'
$ws.Cells.Item(130,3).Value = 'laliga'
$ws.Cells.Item(130,4).Value = '2023-2024'
$ws.Cells.Item(130,5).Value = 45254.875
$ws.Cells.Item(130,6).Value = 'Alaves'
$ws.Cells.Item(130,7).Value = 3
$ws.Cells.Item(130,8).Value = 'Granada CF'
$ws.Cells.Item(130,9).Value = 1
$ws.Cells.Item(130,10).Value = 1.86
$ws.Cells.Item(130,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(130,12).Value = 1.74
$ws.Cells.Item(130,13).Value = '24/11/2023 20:58'
$ws.Cells.Item(130,14).Value = 3.54
$ws.Cells.Item(130,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(130,16).Value = 3.67
$ws.Cells.Item(130,17).Value = '24/11/2023 20:58'
$ws.Cells.Item(130,18).Value = 4.51
$ws.Cells.Item(130,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(130,20).Value = 5.45
$ws.Cells.Item(130,21).Value = '24/11/2023 20:58'
$ws.Cells.Item(130,22).Value = 'https://www.betexplorer.com/football/spain/laliga/alaves-granada-cf/xWExufcb/'

# Row 131
$ws.Range("A130:V130").Copy()
$ws.Range("A131:V131").PasteSpecial(-4122)
$ws.Cells.Item(131,1).Value = 130
$ws.Cells.Item(131,2).Value = 'spain'
$ws.Cells.Item(131,3).Value = 'laliga'
$ws.Cells.Item(131,4).Value = '2023-2024'
$ws.Cells.Item(131,5).Value = 45255.58333333334
$ws.Cells.Item(131,6).Value = 'Rayo Vallecano'
$ws.Cells.Item(131,7).Value = 1
$ws.Cells.Item(131,8).Value = 'Barcelona'
$ws.Cells.Item(131,9).Value = 1
$ws.Cells.Item(131,10).Value = 4.22
$ws.Cells.Item(131,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(131,12).Value = 5
$ws.Cells.Item(131,13).Value = '25/11/2023 13:59'
$ws.Cells.Item(131,14).Value = 3.78
$ws.Cells.Item(131,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(131,16).Value = 4.25
$ws.Cells.Item(131,17).Value = '25/11/2023 13:59'
$ws.Cells.Item(131,18).Value = 1.85
$ws.Cells.Item(131,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(131,20).Value = 1.68
$ws.Cells.Item(131,21).Value = '25/11/2023 13:56'
$ws.Cells.Item(131,22).Value = 'https://www.betexplorer.com/football/spain/laliga/rayo-vallecano-barcelona/2wClxhsH/'

# Row 132
$ws.Range("A131:V131").Copy()
$ws.Range("A132:V132").PasteSpecial(-4122)
$ws.Cells.Item(132,1).Value = 131
$ws.Cells.Item(132,2).Value = 'spain'
$ws.Cells.Item(132,3).Value = 'laliga'
$ws.Cells.Item(132,4).Value = '2023-2024'
$ws.Cells.Item(132,5).Value = 45255.67708333334
$ws.Cells.Item(132,6).Value = 'Valencia'
$ws.Cells.Item(132,7).Value = 0
$ws.Cells.Item(132,8).Value = 'Celta Vigo'
$ws.Cells.Item(132,9).Value = 0
$ws.Cells.Item(132,10).Value = 2
$ws.Cells.Item(132,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(132,12).Value = 2.3
$ws.Cells.Item(132,13).Value = '25/11/2023 16:15'
$ws.Cells.Item(132,14).Value = 3.38
$ws.Cells.Item(132,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(132,16).Value = 3.32
$ws.Cells.Item(132,17).Value = '25/11/2023 16:13'
$ws.Cells.Item(132,18).Value = 4.1
$ws.Cells.Item(132,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(132,20).Value = 3.44
$ws.Cells.Item(132,21).Value = '25/11/2023 16:15'
$ws.Cells.Item(132,22).Value = 'https://www.betexplorer.com/football/spain/laliga/valencia-celta-vigo/ryQIKPBq/'

# Row 133
$ws.Range("A132:V132").Copy()
$ws.Range("A133:V133").PasteSpecial(-4122)
$ws.Cells.Item(133,1).Value = 132
$ws.Cells.Item(133,2).Value = 'spain'
$ws.Cells.Item(133,3).Value = 'laliga'
$ws.Cells.Item(133,4).Value = '2023-2024'
$ws.Cells.Item(133,5).Value = 45255.77083333334
$ws.Cells.Item(133,6).Value = 'Getafe'
$ws.Cells.Item(133,7).Value = 2
$ws.Cells.Item(133,8).Value = 'Almeria'
$ws.Cells.Item(133,9).Value = 1
$ws.Cells.Item(133,10).Value = 1.79
$ws.Cells.Item(133,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(133,12).Value = 1.78
$ws.Cells.Item(133,13).Value = '25/11/2023 18:25'
$ws.Cells.Item(133,14).Value = 3.58
$ws.Cells.Item(133,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(133,16).Value = 3.68
$ws.Cells.Item(133,17).Value = '25/11/2023 18:25'
$ws.Cells.Item(133,18).Value = 4.88
$ws.Cells.Item(133,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(133,20).Value = 5.06
$ws.Cells.Item(133,21).Value = '25/11/2023 18:25'
$ws.Cells.Item(133,22).Value = 'https://www.betexplorer.com/football/spain/laliga/getafe-almeria/GMDtvEC4/'

# Row 134
$ws.Range("A133:V133").Copy()
$ws.Range("A134:V134").PasteSpecial(-4122)
$ws.Cells.Item(134,1).Value = 133
$ws.Cells.Item(134,2).Value = 'spain'
$ws.Cells.Item(134,3).Value = 'laliga'
$ws.Cells.Item(134,4).Value = '2023-2024'
$ws.Cells.Item(134,5).Value = 45255.875
$ws.Cells.Item(134,6).Value = 'Atl. Madrid'
$ws.Cells.Item(134,7).Value = 1
$ws.Cells.Item(134,8).Value = 'Mallorca'
$ws.Cells.Item(134,9).Value = 0
$ws.Cells.Item(134,10).Value = 1.43
$ws.Cells.Item(134,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(134,12).Value = 1.36
$ws.Cells.Item(134,13).Value = '25/11/2023 20:51'
$ws.Cells.Item(134,14).Value = 4.44
$ws.Cells.Item(134,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(134,16).Value = 5.1
$ws.Cells.Item(134,17).Value = '25/11/2023 20:57'
$ws.Cells.Item(134,18).Value = 8.35
$ws.Cells.Item(134,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(134,20).Value = 9.85
$ws.Cells.Item(134,21).Value = '25/11/2023 20:57'
$ws.Cells.Item(134,22).Value = 'https://www.betexplorer.com/football/spain/laliga/atl-madrid-mallorca/QioKMCsU/'

# Row 135
$ws.Range("A134:V134").Copy()
$ws.Range("A135:V135").PasteSpecial(-4122)
$ws.Cells.Item(135,1).Value = 134
$ws.Cells.Item(135,2).Value = 'spain'
$ws.Cells.Item(135,3).Value = 'laliga'
$ws.Cells.Item(135,4).Value = '2023-2024'
$ws.Cells.Item(135,5).Value = 45256.58333333334
$ws.Cells.Item(135,6).Value = 'Villarreal'
$ws.Cells.Item(135,7).Value = 3
$ws.Cells.Item(135,8).Value = 'Osasuna'
$ws.Cells.Item(135,9).Value = 1
$ws.Cells.Item(135,10).Value = 1.76
$ws.Cells.Item(135,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(135,12).Value = 2.2
$ws.Cells.Item(135,13).Value = '26/11/2023 13:59'
$ws.Cells.Item(135,14).Value = 3.86
$ws.Cells.Item(135,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(135,16).Value = 3.46
$ws.Cells.Item(135,17).Value = '26/11/2023 13:58'
$ws.Cells.Item(135,18).Value = 4.69
$ws.Cells.Item(135,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(135,20).Value = 3.52
$ws.Cells.Item(135,21).Value = '26/11/2023 13:59'
$ws.Cells.Item(135,22).Value = 'https://www.betexplorer.com/football/spain/laliga/villarreal-osasuna/dvUMJqRk/'

# Row 136
$ws.Range("A135:V135").Copy()
$ws.Range("A136:V136").PasteSpecial(-4122)
$ws.Cells.Item(136,1).Value = 135
$ws.Cells.Item(136,2).Value = 'spain'
$ws.Cells.Item(136,3).Value = 'laliga'
$ws.Cells.Item(136,4).Value = '2023-2024'
$ws.Cells.Item(136,5).Value = 45256.67708333334
$ws.Cells.Item(136,6).Value = 'Real Sociedad'
$ws.Cells.Item(136,7).Value = 2
$ws.Cells.Item(136,8).Value = 'Sevilla'
$ws.Cells.Item(136,9).Value = 1
$ws.Cells.Item(136,10).Value = 1.77
$ws.Cells.Item(136,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(136,12).Value = 1.63
$ws.Cells.Item(136,13).Value = '26/11/2023 16:06'
$ws.Cells.Item(136,14).Value = 3.6
$ws.Cells.Item(136,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(136,16).Value = 3.95
$ws.Cells.Item(136,17).Value = '26/11/2023 16:14'
$ws.Cells.Item(136,18).Value = 5.01
$ws.Cells.Item(136,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(136,20).Value = 6.09
$ws.Cells.Item(136,21).Value = '26/11/2023 16:14'
$ws.Cells.Item(136,22).Value = 'https://www.betexplorer.com/football/spain/laliga/real-sociedad-sevilla/f7PXuzrh/'

# Row 137
$ws.Range("A136:V136").Copy()
$ws.Range("A137:V137").PasteSpecial(-4122)
$ws.Cells.Item(137,1).Value = 136
$ws.Cells.Item(137,2).Value = 'spain'
$ws.Cells.Item(137,3).Value = 'laliga'
$ws.Cells.Item(137,4).Value = '2023-2024'
$ws.Cells.Item(137,5).Value = 45256.77083333334
$ws.Cells.Item(137,6).Value = 'Cadiz CF'
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 'Real Madrid'
$ws.Cells.Item(137,9).Value = 3
$ws.Cells.Item(137,10).Value = 5.82
$ws.Cells.Item(137,11).Value = '05/11/2023 11:03'
$ws.Cells.Item(137,12).Value = 6.9
$ws.Cells.Item(137,13).Value = '26/11/2023 18:29'
$ws.Cells.Item(137,14).Value = 4.35
$ws.Cells.Item(137,15).Value = '05/11/2023 11:03'
$ws.Cells.Item(137,16).Value = 4.7
$ws.Cells.Item(137,17).Value = '26/11/2023 18:29'
$ws.Cells.Item(137,18).Value = 1.56
$ws.Cells.Item(137,19).Value = '05/11/2023 11:03'
$ws.Cells.Item(137,20).Value = 1.48
$ws.Cells.Item(137,21).Value = '26/11/2023 18:29'
$ws.Cells.Item(137,22).Value = 'https://www.betexplorer.com/football/spain/laliga/cadiz-real-madrid/t4LTtGSo/'

$excel.CutCopyMode = 0